# Weekly update: a new price record (2023-08-03, serial date 45141) is
# reported for "Vega Monumental Concepción - Choclo" and inserted as the
# new row 123, pushing all the existing records (old rows 123-178) down
# by one row (new rows 124-179). Sheet dimension grows from R178 to R179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 123; Excel shifts rows 123:178 down
# to 124:179 and extends the used range automatically.
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new record.
$ws.Range("A123").Value = 11
$ws.Range("B123").Value = "Vega Monumental Concepción"
$ws.Range("C123").Value = "Bíobío"
$ws.Range("D123").Value = 45141
$ws.Range("E123").Value = 8
$ws.Range("F123").Value = 100112024
$ws.Range("G123").Value = "Choclo"
$ws.Range("H123").Value = "Dulce o Americano"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 100
$ws.Range("K123").Value = 40000
$ws.Range("L123").Value = 42000
$ws.Range("M123").Value = 41000
$ws.Range("N123").Value = "$/malla 70 unidades"
$ws.Range("O123").Value = "Región de Arica y Parinacota"
$ws.Range("P123").Value = 586
$ws.Range("Q123").Value = 70
$ws.Range("R123").Value = "Hortaliza"
